# Add a new "2022" data column (K) to the statistics table, mirroring the
# existing "2021" column (J): same formatting/styles, new values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy column J (rows 4-14: header year + all data rows) into column K first,
# so the new column inherits the exact same cell styles/number formats as
# column J (year header style, the "0.0" data style, and the bottom-border
# style on row 14).
$ws.Range("J4:J14").Copy($ws.Range("K4"))

# Now overwrite the copied values with the real 2022 figures.
$ws.Range("K4").Value = 2022

$ws.Range("K5").Value = 1.6
$ws.Range("K6").Value = 0.4
$ws.Range("K7").Value = 0.9
$ws.Range("K8").Value = 0.6
$ws.Range("K9").Value = 2.1
$ws.Range("K10").Value = 0.6
$ws.Range("K11").Value = 0.9
$ws.Range("K12").Value = 2.3
$ws.Range("K13").Value = 4.3
$ws.Range("K14").Value = 0.3

# Match the author's final selection/cursor position.
$ws.Range("L7").Select()
